$d = $word.ActiveDocument

# --------------------------------------------------------------------------
# 1) Template placeholders referenced the data root as "d.data.rollup...."
#    which was a leftover debugging artifact; fix every occurrence back to
#    the real path "d.rollup....". This hits all of the merge-field style
#    placeholders throughout the table (title, dates, portfolio/ project
#    fields, project_health ifEQ blocks, etc).
# --------------------------------------------------------------------------
$d.Content.Find.Execute("data.rollup", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "rollup", 2) | Out-Null

# --------------------------------------------------------------------------
# 2) Remove the leftover sample/test payload that had been pasted below the
#    table (the "{ test: 'TEST', report_date: ..., rollup: { portfolios:
#    [...] } }" block, spanning several paragraphs) -- only a single empty
#    paragraph should remain between the table and the section properties.
#
#    NOTE: the table boundary is located via a Range's Tables collection
#    (Range.Tables.Item) rather than Document.Tables.Item -- the latter
#    leaves the document's paragraph/range cache in a bad state in this
#    host. Likewise the cleanup loop walks every paragraph (instead of
#    `break`-ing out early) since early-exit from this loop was unreliable
#    here too.
# --------------------------------------------------------------------------
$tbl = $d.Content.Tables.Item(1)
$tblEnd = $tbl.Range.End

$paraCount = $d.Paragraphs.Count
for ($i = $paraCount; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Start -lt $tblEnd) {
    } else {
        $p.Range.Delete()
    }
}
